{"js": "const replacements = [\n  [\"2024-04-25 Thursday\", \"2024-04-26 Friday\"],\n  [\"216\u00f78=\", \"728\u00f72=\"],\n  [\"740\u00f78=\", \"268\u00f77=\"],\n  [\"681\u00f78=\", \"568\u00f79=\"],\n  [\"484\u00f79=\", \"503\u00f76=\"],\n  [\"858\u00f73=\", \"645\u00f74=\"],\n  [\"342\u00f78=\", \"237\u00f76=\"],\n  [\"620\u00f73=\", \"580\u00f78=\"],\n  [\"568\u00f74=\", \"870\u00f72=\"],\n  [\"780\u00f78=\", \"820\u00f72=\"],\n  [\"581\u00f78=\", \"943\u00f75=\"],\n  [\"469\u00f72=\", \"955\u00f74=\"],\n  [\"897\u00f78=\", \"591\u00f79=\"],\n  [\"104\u00f76=\", \"457\u00f76=\"],\n  [\"674\u00f76=\", \"608\u00f74=\"],\n  [\"758\u00f73=\", \"148\u00f76=\"],\n  [\"259\u00f72=\", \"408\u00f77=\"],\n  [\"809\u00f76=\", \"756\u00f77=\"],\n  [\"373\u00f72=\", \"314\u00f75=\"],\n  [\"486\u00f73=\", \"456\u00f79=\"],\n  [\"823\u00f72=\", \"841\u00f74=\"],\n  [\"695\u00f78=\", \"326\u00f73=\"],\n  [\"575\u00f77=\", \"138\u00f74=\"],\n  [\"677\u00f78=\", \"995\u00f75=\"],\n  [\"366\u00f77=\", \"489\u00f76=\"],\n  [\"898\u00f76=\", \"891\u00f78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@('2024-04-25 Thursday', '2024-04-26 Friday')\n    ,@('216\u00f78=', '728\u00f72=')\n    ,@('740\u00f78=', '268\u00f77=')\n    ,@('681\u00f78=', '568\u00f79=')\n    ,@('484\u00f79=', '503\u00f76=')\n    ,@('858\u00f73=', '645\u00f74=')\n    ,@('342\u00f78=', '237\u00f76=')\n    ,@('620\u00f73=', '580\u00f78=')\n    ,@('568\u00f74=', '870\u00f72=')\n    ,@('780\u00f78=', '820\u00f72=')\n    ,@('581\u00f78=', '943\u00f75=')\n    ,@('469\u00f72=', '955\u00f74=')\n    ,@('897\u00f78=', '591\u00f79=')\n    ,@('104\u00f76=', '457\u00f76=')\n    ,@('674\u00f76=', '608\u00f74=')\n    ,@('758\u00f73=', '148\u00f76=')\n    ,@('259\u00f72=', '408\u00f77=')\n    ,@('809\u00f76=', '756\u00f77=')\n    ,@('373\u00f72=', '314\u00f75=')\n    ,@('486\u00f73=', '456\u00f79=')\n    ,@('823\u00f72=', '841\u00f74=')\n    ,@('695\u00f78=', '326\u00f73=')\n    ,@('575\u00f77=', '138\u00f74=')\n    ,@('677\u00f78=', '995\u00f75=')\n    ,@('366\u00f77=', '489\u00f76=')\n    ,@('898\u00f76=', '891\u00f78=')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
